# Fill in the previously blank "Function" column (F) for the
# Unsat_Clause_Selector signal rows (89-97) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F89").Value = "control signal handled within the circuit (not controller)"
$ws.Range("F90").Value = "debug clear control signal"
$ws.Range("F91").Value = "debug flag signal"
$ws.Range("F92").Value = "fifo information relayed to determine module function"
$ws.Range("F93").Value = "fifo data input to the unsat clause buffer"
$ws.Range("F94").Value = "random input for clause selection"
$ws.Range("F95").Value = "how many unsat clauses are currently in the buffer"
$ws.Range("F96").Value = "selected clause output"
$ws.Range("F97").Value = "ucb overflow flag"
